# Updated cryptos list on Wed Sep 11 07:42:57 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures and reorders a couple of
# coin rows (Avalanche <-> WrappedliquidstakedEther2.0, Mantle <-> Hedera)
# while leaving the rank column (A) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" strings (e.g. "512.99") must stay TEXT, matching
# the source data. Assigning them straight would make Excel auto-convert
# to a number, so we prefix with an apostrophe to force text, then copy
# the (unstyled) sibling cell's Style back onto the cell to avoid leaving
# a stray "quote prefix" style behind.

$ws.Range("D2").Value = '56.515.84'
$ws.Range("E2").Value = '  -1.26%  '

$ws.Range("D3").Value = '2.331.87'
$ws.Range("E3").Value = '  -1.18%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").Value = '''512.99'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  -1.30%  '

$ws.Range("D6").Value = '''132.19'
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = '  -2.21%  '

$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("E8").Value = '  -0.97%  '

$ws.Range("E9").Value = '  -2.98%  '

$ws.Range("E10").Value = '  -0.47%  '

$ws.Range("D11").Value = '''5.30'
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = '  +1.33%  '

$ws.Range("D12").Value = '''0.339'
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = '  -0.90%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '2.746.85'
$ws.Range("E13").Value = '  -0.93%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''23.56'
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = '  -0.85%  '

$ws.Range("D15").Value = '56.515.73'
$ws.Range("E15").Value = '  -0.96%  '

$ws.Range("E16").Value = '  -1.61%  '

$ws.Range("D17").Value = '2.335.07'
$ws.Range("E17").Value = '  +0.43%  '

$ws.Range("D18").Value = '''10.45'
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = '  -0.65%  '

$ws.Range("D19").Value = '''324.95'
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("E20").Value = '  -2.76%  '

$ws.Range("D21").Value = '''6.71'
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = '  +2.12%  '

$ws.Range("D22").Value = '''0.998'
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").Value = '''61.76'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  +1.15%  '

$ws.Range("D24").Value = '''8.72'
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = '  +11.26%  '

$ws.Range("D25").Value = '''0.163'
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = '  +1.14%  '

$ws.Range("E26").Value = '  -0.36%  '

$ws.Range("E27").Value = '  +4.59%  '

$ws.Range("D28").Value = '''168.00'
$ws.Range("D28").Style = $ws.Range("B28").Style

$ws.Range("D29").Value = '''1.67'
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = '  -0.49%  '

$ws.Range("D30").Value = '0.0₃0720'
$ws.Range("E30").Value = '  -3.53%  '

$ws.Range("D31").Value = '''6.12'
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = '  -1.48%  '

$ws.Range("D32").Value = '''18.37'
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  -0.22%  '

$ws.Range("D35").Value = '''1.26'
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = '  +0.72%  '

$ws.Range("D36").Value = '''3.93'
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = '  -2.44%  '

$ws.Range("D37").Value = '''0.887'
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = '  -4.77%  '

$ws.Range("D38").Value = '''153.95'
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = '  +12.02%  '

$ws.Range("E39").Value = '  +0.93%  '

$ws.Range("E40").Value = '  +1.56%  '

$ws.Range("E41").Value = '  -0.87%  '

$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("D43").Value = '''279.31'
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = '  -0.21%  '

$ws.Range("D44").Value = '''5.05'
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = '  -1.72%  '

$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.559'
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = '  -1.10%  '

$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").Value = '''0.0495'
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = '  -2.11%  '

$ws.Range("D48").Value = '''18.17'
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = '  +4.78%  '

$ws.Range("E49").Value = '  +0.78%  '

$ws.Range("D50").Value = '''0.0214'
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = '  -2.16%  '

$ws.Range("D51").Value = '''17.22'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = '  +1.84%  '
